# Orthograph correction of packages references publications
# - Joins the two-line dropletUtils reference into a single line, using "; " as
#   the separator instead of a newline.
# - Updates the current selection on the "references" sheet.
# - Nudges the column widths (A/C and B) down slightly.
# - Best-effort: also nudges the window's tab-ratio (may not survive save in
#   this host, but set anyway for parity with the host's object model).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Fix the dropletUtils reference text: replace the embedded newline
#        between the two citations with "; " so it reads as one paragraph.
$cell = $ws.Cells.Find("Riesenfeld")
if ($cell -eq $null) {
    $cell = $ws.Range("B8")
}
$text = $cell.Value2
$fixed = $text -replace "`r`n", "; "
$fixed = $fixed -replace "`n", "; "
$cell.Value2 = $fixed

# --- 2) Column widths: shrink col A (and the shared col C..IUV style) and
#        col B slightly, matching the target layout.
$ws.Columns.Item(1).ColumnWidth = 9.5
$ws.Columns.Item(2).ColumnWidth = 321.65
$ws.Columns.Item(3).ColumnWidth = 9.5

# --- 3) Selection: move the active cell/selection to B23.
$null = $ws.Range("B23").Select()

# --- 4) Best-effort: window tab ratio (987/1650).
try {
    $excel.ActiveWindow.TabRatio = 987 / 1650
} catch {
}
